# Updated cryptos list on Thu Sep 12 10:36:32 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.048.12'
$ws.Range('E2').Value = '  +2.77%  '
$ws.Range('D3').Value = '2.334.31'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'544.65"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.41%  '
$ws.Range('D6').Value = "'134.58"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('D9').Value = '2.348.88'
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('E13').Value = '  +7.09%  '
$ws.Range('D14').Value = '2.770.38'
$ws.Range('E14').Value = '  +1.47%  '
$ws.Range('D15').Value = "'23.57"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = '58.026.68'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '2.345.10'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = "'333.80"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').Value = "'6.73"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = "'0.997"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').Value = "'0.169"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.98%  '
$ws.Range('D26').Value = "'8.49"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').Value = "'0.999"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +7.77%  '
$ws.Range('E29').Value = '  +5.46%  '
$ws.Range('D30').Value = "'169.64"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('D31').Value = '0.0₃0732'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').Value = "'6.15"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').Value = "'1.03"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +17.36%  '
$ws.Range('D34').Value = "'18.49"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = "'0.998"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = "'4.19"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.05%  '
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = "'1.64"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.08%  '
$ws.Range('D40').Value = "'39.21"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('D41').Value = "'149.57"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').Value = "'0.379"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = "'3.61"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.15%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').Value = "'285.48"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.59%  '
$ws.Range('D45').Value = "'19.21"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.97%  '
$ws.Range('D46').Value = "'0.0926"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('E47').Value = '  +2.27%  '
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('D50').Value = "'17.56"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.44%  '
$ws.Range('D51').Value = "'0.380"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.33%  '
